# Rename "AddressBook"/"Address Book" related class names to "TravelBuddy"
# throughout the UML sequence diagram on slide 1, per the commit:
# Modify "addressbook", "address book", "person" and "persons" to
# "travelbuddy", "place" or "places"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $candidate = $slide.Shapes.Item($i)
        if ($candidate.Id -eq $id) { return $candidate }
    }
    return $null
}

# --- Shape "Rectangle 62" (id=16): ":Address" / "BookParser" (2 paragraphs)
#     becomes a single paragraph ":" + "TravelBuddyParser", and the shape
#     widens slightly to fit the new label.
$shpParser = Get-ShapeById $s 16
$shpParser.Left = 270.6446456692913
$shpParser.Width = 101.35503937007874

$trParser = $shpParser.TextFrame.TextRange
# Drop the first paragraph ("{:Address}") entirely -- this also removes its
# trailing paragraph mark, leaving the second paragraph's run/pPr/endParaRPr
# intact as the sole remaining paragraph.
$trParser.Paragraphs(1, 1).Delete() | Out-Null
# Re-add the leading ":" as its own run in front of what is now the only
# paragraph ("BookParser").
$trParser.InsertBefore(":") | Out-Null
# Rename the remaining "BookParser" run to "TravelBuddyParser".
$trParser.Characters(2, 10).Text = "TravelBuddyParser"

# --- Shape "TextBox 78" (id=79): "undo" + "AddressBook" + "()"
$shpUndo = Get-ShapeById $s 79
$trUndo = $shpUndo.TextFrame.TextRange
$trUndo.Characters(5, 11).Text = "TravelBuddy"

# --- Shape "Rectangle 62" (id=84): ":" + "VersionedAddressBook"
$shpVersioned = Get-ShapeById $s 84
$trVersioned = $shpVersioned.TextFrame.TextRange
$trVersioned.Characters(2, 20).Text = "VersionedTravelBuddy"

# --- Shape "TextBox 87" (id=88): "resetData" + "(" + "ReadOnlyAddressBook" + ")"
$shpReadOnly = Get-ShapeById $s 88
$trReadOnly = $shpReadOnly.TextFrame.TextRange
$trReadOnly.Characters(11, 19).Text = "ReadOnlyTravelBuddy"
